$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-07-07 Sunday" "2024-07-08 Monday"

Replace-Text "180×2=" "380×2="
Replace-Text "846×9=" "655×2="
Replace-Text "712×8=" "120×6="
Replace-Text "381×3=" "398×2="
Replace-Text "545×6=" "399×6="

Replace-Text "114×4=" "246×6="
Replace-Text "921×9=" "178×8="
Replace-Text "357×2=" "144×5="
Replace-Text "429×5=" "910×5="
Replace-Text "944×3=" "538×7="

Replace-Text "874×6=" "212×7="
Replace-Text "277×9=" "528×2="
Replace-Text "945×4=" "966×8="
Replace-Text "554×5=" "657×5="
Replace-Text "734×2=" "587×4="

Replace-Text "211×2=" "190×6="
Replace-Text "768×9=" "520×8="
Replace-Text "994×5=" "759×3="
Replace-Text "444×8=" "172×4="
Replace-Text "759×8=" "575×6="

Replace-Text "411×3=" "926×6="
Replace-Text "963×4=" "782×2="
Replace-Text "382×7=" "684×4="
Replace-Text "948×2=" "277×7="
Replace-Text "705×2=" "448×8="
